$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Bump the "Size: N test case(s))" label on the summary block (D3) from
#    4 to 5, since we are adding a 5th test case (TC5) below.
# ---------------------------------------------------------------------------
$ws.Range("D3").Value = "Size: 5 test case(s))"

# ---------------------------------------------------------------------------
# 2) Append a new "TC5" test-case block, mirroring the 6-row layout used by
#    the existing TC1..TC4 blocks (Test Case ID / Description / Precondition
#    / header / 2 data rows), two blank rows below the last block (37, 38).
#    Copy formatting only (keep existing shared styles) from the TC1 block
#    (rows 6-11) onto the new rows (39-44), then set the actual values.
# ---------------------------------------------------------------------------
# Merge first, then copy formats onto the (now merged) range -- doing it in
# the opposite order causes the merged cells' style to be rewritten to a
# freshly duplicated style instead of reusing the pasted-in shared style.
$ws.Range("B40:D40").Merge()
$ws.Range("B41:F41").Merge()

$ws.Range("A6:F11").Copy()
$ws.Range("A39").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 39: Test Case ID row
$ws.Range("A39").Value = "Test Case ID: "
$ws.Range("B39").Value = "TC5"
$ws.Range("C39").Value = "Priority (low,medium,high: "
$ws.Range("E39").Value = "Executed by:"

# Row 40: Description row
$ws.Range("A40").Value = "Description: "
$ws.Range("E40").Value = "Execution Date: "

# Row 41: Precondition row
$ws.Range("A41").Value = "Precondition: "
$ws.Range("B41").Value = "O usuario devidamente autenticado e na tela inicial do sistema"

# Row 42: table header row
$ws.Range("A42").Value = "#"
$ws.Range("B42").Value = "Steps"
$ws.Range("C42").Value = "Test Data"
$ws.Range("D42").Value = "Expected Results"
$ws.Range("E42").Value = "Execution Status (pass/fail/blocked)"
$ws.Range("F42").Value = "Actual Result"

# Row 43: first test step
$ws.Range("A43").Value = 1
$ws.Range("B43").Value = "Chefe Clica para listar todas as solicitações de diárias relacionadas à sua competência."
$ws.Range("D43").Value = "SYSTEM Recupera os registros de solicitações e os exibe (em ordem decrescente pelo número da diária) em tela para o usuário."

# Row 44: second test step (the new scenario for this test case)
$ws.Range("A44").Value = 2
$ws.Range("B44").Value = "Chefe Clica para ordenar pelo nome do servidor."
$ws.Range("D44").Value = "SYSTEM Visualiza os registros de solicitações de diária ordenado pelo nome do servidor."
